# Auto-generated edit script: reverts accidental inline-added "contributor" image
# references and restores the intended data (see commit message: "don't add inline contributors").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SchemaOrganization")
$ws.Range("B3").Value = "http://example.com/organization3:Image0"
$ws.Range("B4").Value = "http://example.com/organization5:Image1"

$ws = $wb.Worksheets.Item("FoafPerson")
$ws.Range("E3").Value = "http://example.com/person2:Image0"

$ws = $wb.Worksheets.Item("SchemaPerson")
$ws.Range("D2").Value = "http://example.com/person1:Image0"
$ws.Range("D4").Value = "http://example.com/person5:Image0"

$ws = $wb.Worksheets.Item("SchemaCreativeWork")
$ws.Range("E2").Value = "https://images.metmuseum.org/CRDImages/ep/original/LC-EP_1993_132_suppl_CH-002.jpg"

$ws = $wb.Worksheets.Item("RdfProperty")
$ws.Range("C2").Value = "dcterms:description:Image0"
$ws.Range("C3").Value = "dcterms:extent:Image0"
$ws.Range("C4").Value = "dcterms:language:Image1"
$ws.Range("C6").Value = "dcterms:publisher:Image0"
$ws.Range("C8").Value = "dcterms:spatial:Image0"
$ws.Range("C11").Value = "dcterms:type:Image0"

$ws = $wb.Worksheets.Item("SchemaProperty")
$ws.Range("C3").Value = "schema:description:Image0"
$ws.Range("C5").Value = "schema:spatial:Image1"

$ws = $wb.Worksheets.Item("FoafOrganization")
$ws.Range("C3").Value = "http://example.com/organization2:Image0"
$ws.Range("C4").Value = "http://example.com/organization4:Image0"

$ws = $wb.Worksheets.Item("SkosConcept")
$ws.Range("B3").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:1:Image1"
$ws.Range("B5").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:3:Image0"
$ws.Range("B6").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:4:Image0"
$ws.Range("B8").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:6:Image0"
$ws.Range("B9").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:7:Image1"
$ws.Range("B11").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:9:Image1"
$ws.Range("B13").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:11:Image1"
$ws.Range("B14").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:12:Image0"
$ws.Range("B16").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:14:Image0"
$ws.Range("B20").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:18:Image1"
$ws.Range("B22").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:20:Image1"
$ws.Range("B23").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:21:Image0"
$ws.Range("B24").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:22:Image1"
$ws.Range("B28").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:26:Image1"
$ws.Range("B34").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:32:Image0"
$ws.Range("B36").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:34:Image0"
$ws.Range("B37").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:35:Image0"
$ws.Range("B39").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:37:Image1"
$ws.Range("B40").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:38:Image1"
$ws.Range("B41").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:39:Image0"
$ws.Range("B42").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:40:Image1"
$ws.Range("B44").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:42:Image1"
$ws.Range("B45").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:43:Image1"
$ws.Range("B46").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:44:Image1"
$ws.Range("B48").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:46:Image0"
$ws.Range("B51").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:49:Image1"
$ws.Range("B55").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:53:Image1"
$ws.Range("B56").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:54:Image0"
$ws.Range("B60").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:58:Image1"
$ws.Range("B64").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:62:Image1"
$ws.Range("B68").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:66:Image0"
$ws.Range("B72").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:70:Image1"
$ws.Range("B73").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:71:Image1"
$ws.Range("B74").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:72:Image1"
$ws.Range("B76").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:74:Image0"
$ws.Range("B77").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:75:Image0"
$ws.Range("B81").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:79:Image1"

$ws = $wb.Worksheets.Item("SchemaDefinedTerm")
$ws.Range("B2").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:80:Image1"
$ws.Range("B4").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:82:Image0"
$ws.Range("B5").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:83:Image0"
$ws.Range("B7").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:85:Image0"
$ws.Range("B9").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:87:Image0"
$ws.Range("B15").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:93:Image0"
$ws.Range("B19").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:97:Image1"
$ws.Range("B21").Value = "urn:paradicms:etl:pipeline:synthetic_data:concept:99:Image0"

$ws = $wb.Worksheets.Item("CreativeCommonsLicense")
$ws.Range("A3").Value = "http://creativecommons.org/licenses/by/4.0/"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "4.0"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "by"
$ws.Range("J3").Value = "Attribution 4.0 International"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("A4").Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "1.0"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "nc"
$ws.Range("L4").Value = "NonCommercial 1.0 Generic"
$ws.Range("A5").Value = "http://creativecommons.org/publicdomain/mark/1.0/"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "1.0"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "mark"
$ws.Range("H5").Value = "Public Domain Mark 1.0"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""

$ws = $wb.Worksheets.Item("RightsStatementsDotOrgRightsStatement")
$ws.Range("A2").Value = "http://rightsstatements.org/vocab/InC/1.0/"
$ws.Range("B2").Value = "This Item is protected by copyright and/or related rights.`n`n  You are free to use this Item in any way that is permitted by the copyright and related rights legislation that applies to your use.`n`n  For other uses you need to obtain permission from the rights-holder(s)."
$ws.Range("C2").Value = "This Rights Statement indicates that the Item labeled with this Rights Statement is in copyright."
$ws.Range("D2").Value = "InC"
$ws.Range("E2").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$ws.Range("F2").Value = "In Copyright"
$ws.Range("G2").Value = "This Rights Statement can be used for an Item that is in copyright. Using this statement implies that the organization making this Item available has determined that the Item is in copyright and either is the rights-holder, has obtained permission from the rights-holder(s) to make their Work(s) available, or makes the Item available under an exception or limitation to copyright (including Fair Use) that entitles it to make the Item available."
$ws.Range("A3").Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Range("B3").Value = "This Item is protected by copyright and/or related rights.`n`n  You are free to use this Item in any way that is permitted by the copyright and related rights legislation that applies to your use. In addition, no permission is required from the rights-holder(s) for educational uses.`n`n  For other uses, you need to obtain permission from the rights-holder(s)."
$ws.Range("C3").Value = "This Rights Statement indicates that the Item labeled with this Rights Statement is in copyright but that educational use is allowed without the need to obtain additional permission."
$ws.Range("D3").Value = "InC-EDU"
$ws.Range("E3").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$ws.Range("F3").Value = "In Copyright - Educational Use Permitted"
$ws.Range("G3").Value = "This Rights Statement can be used only for copyrighted Items for which the organization making the Item available is the rights-holder or has been explicitly authorized by the rights-holder(s) to allow third parties to use their Work(s) for educational purposes without first obtaining permission."
$ws.Range("A4").Value = "http://rightsstatements.org/vocab/NoC-US/1.0/"
$ws.Range("B4").Value = "The organization that has made the Item available believes that the Item is in the Public Domain under the laws of the United States, but a determination was not made as to its copyright status under the copyright laws of other countries. The Item may not be in the Public Domain under the laws of other countries.`n`n  Please refer to the organization that has made the Item available for more information."
$ws.Range("C4").Value = "This Rights Statement indicates that the Item is in the Public Domain under the laws of the United States, but that a determination was not made as to its copyright status under the copyright laws of other countries."
$ws.Range("D4").Value = "NoC-US"
$ws.Range("E4").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$ws.Range("F4").Value = "No Copyright - United States"
$ws.Range("G4").Value = "This Rights Statement should be used for Items for which the organization that intends to make the Item available has determined are free of copyright under the laws of the United States. This Rights Statement should not be used for Orphan Works (which are assumed to be in-copyright) or for Items where the organization that intends to make the Item available has not undertaken an effort to ascertain the copyright status of the underlying Work."
$ws.Range("A5").Value = "https://rightsstatements.org/vocab/NoC-US/1.0/"
$ws.Range("B5").Value = "The organization that has made the Item available believes that the Item is in the Public Domain under the laws of the United States, but a determination was not made as to its copyright status under the copyright laws of other countries. The Item may not be in the Public Domain under the laws of other countries.`n`n  Please refer to the organization that has made the Item available for more information."
$ws.Range("C5").Value = "This Rights Statement indicates that the Item is in the Public Domain under the laws of the United States, but that a determination was not made as to its copyright status under the copyright laws of other countries."
$ws.Range("D5").Value = "NoC-US"
$ws.Range("E5").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$ws.Range("F5").Value = "No Copyright - United States"
$ws.Range("G5").Value = "This Rights Statement should be used for Items for which the organization that intends to make the Item available has determined are free of copyright under the laws of the United States. This Rights Statement should not be used for Orphan Works (which are assumed to be in-copyright) or for Items where the organization that intends to make the Item available has not undertaken an effort to ascertain the copyright status of the underlying Work."

